# Adds one new weekly price-report block (2 rows: Primera / Segunda quality)
# for "Pepino ensalada" just above the existing row 315, shifting the rest
# of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 315-316 (pushes old rows 315..332 down to 317..334)
$ws.Range("A315:A316").EntireRow.Insert()

# New row 315 - Calidad "Primera"
$ws.Cells.Item(315, 1).Value = 1
$ws.Cells.Item(315, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(315, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(315, 4).Value = 44753
$ws.Cells.Item(315, 5).Value = 15
$ws.Cells.Item(315, 6).Value = 100112043
$ws.Cells.Item(315, 7).Value = "Pepino ensalada"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 120
$ws.Cells.Item(315, 11).Value = 14000
$ws.Cells.Item(315, 12).Value = 15000
$ws.Cells.Item(315, 13).Value = 14500
$ws.Cells.Item(315, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(315, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(315, 16).Value = 207
$ws.Cells.Item(315, 17).Value = 70
$ws.Cells.Item(315, 18).Value = "Hortaliza"

# New row 316 - Calidad "Segunda"
$ws.Cells.Item(316, 1).Value = 1
$ws.Cells.Item(316, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(316, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(316, 4).Value = 44753
$ws.Cells.Item(316, 5).Value = 15
$ws.Cells.Item(316, 6).Value = 100112043
$ws.Cells.Item(316, 7).Value = "Pepino ensalada"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Segunda"
$ws.Cells.Item(316, 10).Value = 120
$ws.Cells.Item(316, 11).Value = 11000
$ws.Cells.Item(316, 12).Value = 12000
$ws.Cells.Item(316, 13).Value = 11500
$ws.Cells.Item(316, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(316, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(316, 16).Value = 115
$ws.Cells.Item(316, 17).Value = 100
$ws.Cells.Item(316, 18).Value = "Hortaliza"
